# Refresh the "Narastające" (cumulative) AD:AH data and the active cell
# selection on the sole worksheet, matching the upstream data refresh that
# produced the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where only AD changes (previously a blank "-" placeholder, now a
# real number) because B/C/E/F/etc. values on these rows don't extend that
# far into the month.
$ws.Range("AD2").Value = 1918.826086956522
$ws.Range("AD3").Value = 993.26086956521749
$ws.Range("AD4").Value = 1727.0869565217392
$ws.Range("AD8").Value = 2439.2147806004618
$ws.Range("AD20").Value = 2677.6956521739135
$ws.Range("AD21").Value = 2487.130434782609
$ws.Range("AD22").Value = 2471.5813953488373
$ws.Range("AD26").Value = 2544.3913043478265
$ws.Range("AD27").Value = 1377.2207084468666
$ws.Range("AD28").Value = 2113.7647058823532
$ws.Range("AD32").Value = 1358.3259911894274
$ws.Range("AD33").Value = 1290.7826086956522
$ws.Range("AD34").Value = 1770.7826086956522
$ws.Range("AD38").Value = 570.41474654377885
$ws.Range("AD39").Value = 1249.0434782608697
$ws.Range("AD44").Value = 6876.3380281690133
$ws.Range("AD45").Value = 4978.7264150943402
$ws.Range("AD46").Value = 3610.864864864865
$ws.Range("AD50").Value = 644.73913043478262
$ws.Range("AD51").Value = 758.75000000000011
$ws.Range("AD57").Value = 2290.3174603174602
$ws.Range("AD58").Value = 3766.5198237885465
$ws.Range("AD62").Value = 2262.5720620842571
$ws.Range("AD63").Value = 4338.913043478261
$ws.Range("AD64").Value = 2483.4782608695655

# Rows where the trailing AD:AH block (the forward-filled average for the
# remaining days of the month) was recalculated to a new value.
$ws.Range("AD5:AH5").Value = 2328.3597883597881
$ws.Range("AD6:AH6").Value = 2321.5230842455971
$ws.Range("AD7:AH7").Value = 2556.4864864864867
$ws.Range("AD23:AH23").Value = 3510.0036027380816
$ws.Range("AD24:AH24").Value = 3022.1283290257329
$ws.Range("AD25:AH25").Value = 3363.4242493696997
$ws.Range("AD29:AH29").Value = 2161.9519484998277
$ws.Range("AD30:AH30").Value = 1797.6758476901127
$ws.Range("AD31:AH31").Value = 2612.507950931395
$ws.Range("AD35:AH35").Value = 1300.1344537815128
$ws.Range("AD36:AH36").Value = 1304.172077922078
$ws.Range("AD37:AH37").Value = 1215.279589934762
$ws.Range("AD41:AH41").Value = 1643.4579439252338
$ws.Range("AD42:AH42").Value = 1758.5545023696682
$ws.Range("AD47:AH47").Value = 3948.3606345930598
$ws.Range("AD48:AH48").Value = 3634.4825864436398
$ws.Range("AD49:AH49").Value = 4031.1685912240187
$ws.Range("AD53:AH53").Value = 5051.9765548681216
$ws.Range("AD54:AH54").Value = 4652.6335877862593
$ws.Range("AD60:AH60").Value = 2329.4794264339157
$ws.Range("AD61:AH61").Value = 2681.0063872877395
$ws.Range("AD65:AH65").Value = 2918.8154375238823
$ws.Range("AD66:AH66").Value = 2468.0395794681513
$ws.Range("AD67:AH67").Value = 2748.6380597014922

# Row 11's forward-fill starts a column earlier (from F, not AD) and was
# also recalculated to a new value.
$ws.Range("F11:AH11").Value = 1933.3094331167379

# Update the active cell/selection left in the sheet after the refresh.
$ws.Range("F10").Select()
